$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Paciente: nombre / no. expediente
$ws.Range("A6").Value = "MONTERROSO  LÓPEZ  CRISTIAN  JOSUÉ"
$ws.Range("G6").Value = "/201773491"

# Fecha de nacimiento / edad / lugar de nacimiento
# (A9 looks like a date and D9 looks like a plain number to the
# auto-detector, so a leading apostrophe keeps them stored as literal
# text, same as the original cells.)
$ws.Range("A9").Value = "'1993-02-23"
$ws.Range("D9").Value = "'24"
$ws.Range("E9").Value = "CHINAUTLA"

# Ocupacion / nacionalidad / documento de identificacion
$ws.Range("C11").Value = "VARIOS"
$ws.Range("E11").Value = "GUATEMALTECO"
$ws.Range("G11").Value = "'2201016710106"

# Datos de emergencia: nombre / parentesco / direccion / telefono
$ws.Range("A13").Value = "ANDREA MORENO"
$ws.Range("D13").Value = "ENCARGADA"
$ws.Range("E13").Value = "34 AV 12-16 Z. 5 EL EDEN"
$ws.Range("G13").Value = "'49803299"

# Hora y fecha de asistencia medica
$ws.Range("D14").Value = "Hora: 15:45:35"
$ws.Range("A15").Value = "20/11/2017"
